$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture column M's width so the newly inserted column N can match it (11 chars, no bestFit)
$mWidth = $ws.Columns("M").ColumnWidth

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make this sheet the active sheet/tab, and set its selection
$ws.Activate() | Out-Null
$ws.Range("J15").Select() | Out-Null
